$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format so numeric-looking strings (e.g. "32.50", "0.940")
# keep their exact original formatting instead of being parsed as numbers.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.411.15'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -4.43%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.988.71'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -5.80%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '572.14'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -3.42%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '125.13'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -7.28%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '2.981.21'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -6.01%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -5.87%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.06'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -3.74%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.435'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -4.31%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000223'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -6.06%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.50'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -7.12%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.22%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.484.34'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -5.70%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '60.427.66'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -4.24%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.989.66'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -5.63%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.14'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -6.83%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '430.45'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -6.72%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.99'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -6.60%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.658'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -6.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.13'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -6.80%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.90'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -4.45%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '78.85'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -5.53%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.13%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -6.49%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.10'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -8.51%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.88'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -7.86%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '25.14'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -7.58%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.98'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -11.64%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0925'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -10.17%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.25'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -5.12%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.940'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -9.38%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.53'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -5.47%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '49.38'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -3.98%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0₃0647'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -8.77%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -8.45%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -3.79%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -5.73%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '372.26'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -8.29%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.652.51'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -5.82%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.38'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -8.53%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -7.75%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '118.77'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -4.23%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.94'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -8.79%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -5.19%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '23.31'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -8.30%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '31.54'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -7.37%  '
